$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (TungDo) - history strings updated
$ws.Range("F5").Value = ";0;0;0;0;0;0"
$ws.Range("G5").Value = ";14;43;33;3;3;0"
$ws.Range("H5").Value = ";-100;-400.0;-400.0;-215.0;-215.0;-100"

# Row 6 (anhlavodich) - balance/wins/losses + history strings updated
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = ";1;1;0;0;0;1;0;0"
$ws.Range("G6").Value = ";22;23;43;22;41;30;30;24"
$ws.Range("H6").Value = ";+100;+100;-2200;-200;-1000;+100;-250;-100"

# Row 9 (new user: accmoii)
$ws.Range("A9").Value = "accmoii"
$ws.Range("B9").Value = 150
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = ";0;0;1"
$ws.Range("G9").Value = ";14;33;13"
$ws.Range("H9").Value = ";-120.0;-120;+50.0"

# Stray formatted cell left over at E15 (border only, no value)
$ws.Range("E15").Style = "Normal"
$ws.Range("E15").Borders.Item(7).Weight = 1
$ws.Range("E15").Borders.Item(7).LineStyle = 1

# Leave selection on H9, matching the last cell the author edited
$ws.Range("H9").Select()
